$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.48
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 1.74
$ws.Range("G4").Value = 3.5
$ws.Range("H4").Value = 2.78
$ws.Range("K4").Value = 2.92
$ws.Range("Q4").Value = 3.6
$ws.Range("S4").Value = 8.2
$ws.Range("T4").Value = 2.64
$ws.Range("U4").Value = 1.45
$ws.Range("X4").Value = 6
$ws.Range("AD4").Value = 17
$ws.Range("AH4").Value = 1000
$ws.Range("F5").Value = 2.14
$ws.Range("J5").Value = 2.96
$ws.Range("K5").Value = 3.35
$ws.Range("P5").Value = 1.51
$ws.Range("U5").Value = 1.64
$ws.Range("V5").Value = 1.29
$ws.Range("X5").Value = 9.8
$ws.Range("Y5").Value = 980
$ws.Range("AB5").Value = 980
$ws.Range("AD5").Value = 980
$ws.Range("AG5").Value = 14.5
$ws.Range("F6").Value = 2.06
$ws.Range("G6").Value = 2.16
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 5.1
$ws.Range("J6").Value = 3.05
$ws.Range("N6").Value = 2.68
$ws.Range("O6").Value = 1.52
$ws.Range("P6").Value = 1.56
$ws.Range("Q6").Value = 2.28
$ws.Range("S6").Value = 4.4
$ws.Range("W6").Value = 1.86
$ws.Range("Z6").Value = 1000
$ws.Range("AI6").Value = 120
$ws.Range("F7").Value = 2.32
$ws.Range("G7").Value = 2.48
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 3.55
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 3.5
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 2.48
$ws.Range("P7").Value = 1.5
$ws.Range("Q7").Value = 2.72
$ws.Range("S7").Value = 5.2
$ws.Range("V7").Value = 1.39
$ws.Range("W7").Value = 1.67
$ws.Range("X7").Value = 8.6
$ws.Range("Y7").Value = 9.2
$ws.Range("Z7").Value = 23
$ws.Range("AA7").Value = 85
$ws.Range("AB7").Value = 7.6
$ws.Range("AC7").Value = 8.2
$ws.Range("AD7").Value = 17
$ws.Range("AE7").Value = 65
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 27
$ws.Range("AI7").Value = 95
$ws.Range("AJ7").Value = 38
$ws.Range("AL7").Value = 75
$ws.Range("AM7").Value = 240
$ws.Range("AN7").Value = 980
$ws.Range("AO7").Value = 100
